$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H10").Value = 9363.637000000001
$ws.Range("I10").Value = 7333.3335
$ws.Range("J10").Value = 18500
$ws.Range("K10").Value = 7333.3335
$ws.Range("L10").Value = 18500
$ws.Range("M10").Value = -7040.3335
$ws.Range("N10").Value = -19086

$ws.Range("H101").Value = 604.2857
$ws.Range("I101").Value = 457.5
$ws.Range("J101").Value = 800
$ws.Range("K101").Value = 1372.5
$ws.Range("L101").Value = 2400
$ws.Range("M101").Value = 249.5
$ws.Range("N101").Value = -5644

$ws.Range("H103").Value = 712.3
$ws.Range("I103").Value = 437.2
$ws.Range("K103").Value = 1311.6
$ws.Range("M103").Value = -725.5999999999999

$ws.Range("H133").Value = 12662.308
$ws.Range("J133").Value = 12662.308
$ws.Range("L133").Value = 12662.308
$ws.Range("N133").Value = -22782.308

$ws.Range("H135").Value = 939.80646
$ws.Range("I135").Value = 924.9167
$ws.Range("J135").Value = 990.8570999999999
$ws.Range("K135").Value = 8324.2503
$ws.Range("L135").Value = 8917.713899999999
$ws.Range("M135").Value = -5789.2503
$ws.Range("N135").Value = -13987.7139

$ws.Range("H138").Value = 4105437.2
$ws.Range("I138").Value = 1907340.9
$ws.Range("J138").Value = 4764866
$ws.Range("K138").Value = 5722022.699999999
$ws.Range("L138").Value = 14294598
$ws.Range("M138").Value = -5716882.699999999
$ws.Range("N138").Value = -14304878

$ws.Range("H141").Value = 1411.8438
$ws.Range("I141").Value = 1296.742
$ws.Range("J141").Value = 4980
$ws.Range("K141").Value = 3890.226
$ws.Range("L141").Value = 14940
$ws.Range("M141").Value = 1289.774
$ws.Range("N141").Value = -25300

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 3608.2
$ws.Range("I61").Value = 2858.2307
$ws.Range("J61").Value = 4420.6665
$ws.Range("K61").Value = 2858.2307
$ws.Range("L61").Value = 4420.6665
$ws.Range("M61").Value = -2646.2307
$ws.Range("N61").Value = -4844.6665

$ws.Range("H74").Value = 7528.8286
$ws.Range("I74").Value = 1272.7241
$ws.Range("K74").Value = 1272.7241
$ws.Range("M74").Value = -398.7240999999999

$ws.Range("H77").Value = 7528.8286
$ws.Range("I77").Value = 1272.7241
$ws.Range("K77").Value = 6363.620499999999
$ws.Range("M77").Value = -1995.620499999999

$ws.Range("H132").Value = 4511.913
$ws.Range("I132").Value = 4490.0557
$ws.Range("J132").Value = 4590.6
$ws.Range("K132").Value = 13470.1671
$ws.Range("L132").Value = 13771.8
$ws.Range("M132").Value = -10940.1671
$ws.Range("N132").Value = -18831.8

$ws.Range("H136").Value = 3608.2
$ws.Range("I136").Value = 2858.2307
$ws.Range("J136").Value = 4420.6665
$ws.Range("K136").Value = 8574.6921
$ws.Range("L136").Value = 13261.9995
$ws.Range("M136").Value = -6024.6921
$ws.Range("N136").Value = -18361.9995

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H59").Value = 40575
$ws.Range("J59").Value = 47433.332
$ws.Range("L59").Value = 47433.332
$ws.Range("N59").Value = -49127.332

$ws.Range("H134").Value = 3827.8386
$ws.Range("I134").Value = 2285.353
$ws.Range("K134").Value = 6856.059
$ws.Range("M134").Value = -4321.059

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H10").Value = 569.4
$ws.Range("I10").Value = 569.4
$ws.Range("K10").Value = 569.4
$ws.Range("M10").Value = -430.4

$ws.Range("H13").Value = 205400
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 256250
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 256250
$ws.Range("M13").Value = -1861
$ws.Range("N13").Value = -256528

$ws.Range("H31").Value = 1704.4242
$ws.Range("I31").Value = 1029.6364
$ws.Range("J31").Value = 3054
$ws.Range("K31").Value = 1029.6364
$ws.Range("L31").Value = 3054
$ws.Range("M31").Value = -734.6364000000001
$ws.Range("N31").Value = -3644

$ws.Range("H34").Value = 1704.4242
$ws.Range("I34").Value = 1029.6364
$ws.Range("J34").Value = 3054
$ws.Range("K34").Value = 1029.6364
$ws.Range("L34").Value = 3054
$ws.Range("M34").Value = -827.6364000000001
$ws.Range("N34").Value = -3458

$ws.Range("H58").Value = 2481.2
$ws.Range("I58").Value = 1616.3846
$ws.Range("K58").Value = 1616.3846
$ws.Range("M58").Value = -1413.3846

$ws.Range("H122").Value = 1509.5
$ws.Range("I122").Value = 1073.7142
$ws.Range("J122").Value = 2017.9166
$ws.Range("K122").Value = 3221.1426
$ws.Range("L122").Value = 6053.7498
$ws.Range("M122").Value = -771.1425999999997
$ws.Range("N122").Value = -10953.7498

$ws.Range("H125").Value = 54900
$ws.Range("J125").Value = 54900
$ws.Range("L125").Value = 54900
$ws.Range("N125").Value = -59820

$ws.Range("H132").Value = 2403.0833
$ws.Range("I132").Value = 1995.7
$ws.Range("J132").Value = 4440
$ws.Range("K132").Value = 5987.1
$ws.Range("L132").Value = 13320
$ws.Range("M132").Value = -3457.1
$ws.Range("N132").Value = -18380

$ws.Range("H134").Value = 3714.4211
$ws.Range("I134").Value = 1942.4615
$ws.Range("J134").Value = 7553.6665
$ws.Range("K134").Value = 5827.3845
$ws.Range("L134").Value = 22660.9995
$ws.Range("M134").Value = -3292.3845
$ws.Range("N134").Value = -27730.9995

$ws.Range("H136").Value = 2481.2
$ws.Range("I136").Value = 1616.3846
$ws.Range("K136").Value = 4849.1538
$ws.Range("M136").Value = -2299.1538

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H4").Value = 20146
$ws.Range("I4").Value = 123.166664
$ws.Range("J4").Value = 68200.8
$ws.Range("K4").Value = 369.499992
$ws.Range("L4").Value = 204602.4
$ws.Range("M4").Value = -257.499992
$ws.Range("N4").Value = -204826.4

$ws.Range("H51").Value = 866.6667
$ws.Range("I51").Value = 860
$ws.Range("J51").Value = 900
$ws.Range("K51").Value = 2580
$ws.Range("L51").Value = 2700
$ws.Range("M51").Value = -2120
$ws.Range("N51").Value = -3620

$ws.Range("H113").Value = 9616102
$ws.Range("I113").Value = 508.23077
$ws.Range("J113").Value = 19231696
$ws.Range("K113").Value = 1524.69231
$ws.Range("L113").Value = 57695088
$ws.Range("M113").Value = 645.3076900000001
$ws.Range("N113").Value = -57699428

$ws.Range("H131").Value = 1447.614
$ws.Range("J131").Value = 1677.8511
$ws.Range("L131").Value = 5033.5533
$ws.Range("N131").Value = -15113.5533

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888

$ws.Range("H9").Value = 1003.5
$ws.Range("I9").Value = 1003.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1003.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -833.5
$ws.Range("N9").Value = ""

$ws.Range("H80").Value = 2478.6316
$ws.Range("I80").Value = 2492
$ws.Range("J80").Value = 2455.7144
$ws.Range("K80").Value = 2492
$ws.Range("L80").Value = 2455.7144
$ws.Range("M80").Value = -1494
$ws.Range("N80").Value = -4451.7144

$ws.Range("H83").Value = 2478.6316
$ws.Range("I83").Value = 2492
$ws.Range("J83").Value = 2455.7144
$ws.Range("K83").Value = 12460
$ws.Range("L83").Value = 12278.572
$ws.Range("M83").Value = -7468
$ws.Range("N83").Value = -22262.572

$ws.Range("H104").Value = 55555.555
$ws.Range("J104").Value = 55555.555
$ws.Range("L104").Value = 55555.555
$ws.Range("N104").Value = -62543.555

$ws.Range("H122").Value = 1588733.4
$ws.Range("I122").Value = 2778644
$ws.Range("J122").Value = 2186
$ws.Range("K122").Value = 8335932
$ws.Range("L122").Value = 6558
$ws.Range("M122").Value = -8333482
$ws.Range("N122").Value = -11458

$ws.Range("H126").Value = 2471.7856
$ws.Range("I126").Value = 1716.5834
$ws.Range("J126").Value = 3038.1875
$ws.Range("K126").Value = 5149.7502
$ws.Range("L126").Value = 9114.5625
$ws.Range("M126").Value = -2679.7502
$ws.Range("N126").Value = -14054.5625

$ws.Range("H132").Value = 5435.5
$ws.Range("I132").Value = 5751.9165
$ws.Range("K132").Value = 17255.7495
$ws.Range("M132").Value = -14725.7495

$ws.Range("H138").Value = 62500
$ws.Range("J138").Value = 62500
$ws.Range("L138").Value = 62500
$ws.Range("N138").Value = -72780

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H9").Value = 1399.75
$ws.Range("I9").Value = 933
$ws.Range("J9").Value = 2800
$ws.Range("K9").Value = 933
$ws.Range("L9").Value = 2800
$ws.Range("M9").Value = -709
$ws.Range("N9").Value = -3248

$ws.Range("H125").Value = 50500
$ws.Range("J125").Value = 50500
$ws.Range("L125").Value = 50500
$ws.Range("N125").Value = -60340

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 52867
$ws.Range("I122").Value = 73471.92999999999
$ws.Range("J122").Value = 4788.8335
$ws.Range("K122").Value = 220415.79
$ws.Range("L122").Value = 14366.5005
$ws.Range("M122").Value = -217965.79
$ws.Range("N122").Value = -19266.5005

$ws.Range("H132").Value = 7813891
$ws.Range("I132").Value = 11112304
$ws.Range("J132").Value = 1859.2106
$ws.Range("K132").Value = 33336912
$ws.Range("L132").Value = 5577.6318
$ws.Range("M132").Value = -33334382
$ws.Range("N132").Value = -10637.6318
